$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Il34"
$ws.Range("C2").Value = "Ptprz1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.024860333333333
$ws.Range("H2").Value = 6.074581
$ws.Range("I2").Value = 0.06209502815843435
$ws.Range("J2").Value = 0.06209502815843435
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.06163433333333333
$ws.Range("N2").Value = 0.184903
$ws.Range("O2").Value = 0.006690894379667537
$ws.Range("P2").Value = 0.006690894379667537
$ws.Range("Q2").Value = 0.1248009167381111
$ws.Range("R2").Value = 1.123208250643
$ws.Range("S2").Value = 0.0004154712749105658
$ws.Range("T2").Value = 0.0004154712749105659
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Il34"
$ws.Range("C3").Value = "Ptprz1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.024860333333333
$ws.Range("H3").Value = 6.074581
$ws.Range("I3").Value = 0.06209502815843435
$ws.Range("J3").Value = 0.06209502815843435
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.1030763333333333
$ws.Range("N3").Value = 0.309229
$ws.Range("O3").Value = 0.01118975126488057
$ws.Range("P3").Value = 0.01118975126488057
$ws.Range("Q3").Value = 0.2087151786721111
$ws.Range("R3").Value = 1.878436608049
$ws.Range("S3").Value = 0.0006948279198786356
$ws.Range("T3").Value = 0.0006948279198786356
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Il34"
$ws.Range("C4").Value = "Ptprz1"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.024860333333333
$ws.Range("H4").Value = 6.074581
$ws.Range("I4").Value = 0.06209502815843435
$ws.Range("J4").Value = 0.06209502815843435
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 9.046962666666667
$ws.Range("N4").Value = 27.140888
$ws.Range("O4").Value = 0.9821193543554519
$ws.Range("P4").Value = 0.9821193543554518
$ws.Range("Q4").Value = 18.31883584088089
$ws.Range("R4").Value = 164.869522567928
$ws.Range("S4").Value = 0.06098472896364514
$ws.Range("T4").Value = 0.06098472896364514
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Il34"
$ws.Range("C5").Value = "Ptprz1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 2.952820666666666
$ws.Range("H5").Value = 8.858461999999999
$ws.Range("I5").Value = 0.09055216274676732
$ws.Range("J5").Value = 0.09055216274676733
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.06163433333333333
$ws.Range("N5").Value = 0.184903
$ws.Range("O5").Value = 0.006690894379667537
$ws.Range("P5").Value = 0.006690894379667537
$ws.Range("Q5").Value = 0.1819951332428889
$ws.Range("R5").Value = 1.637956199186
$ws.Range("S5").Value = 0.0006058749567890856
$ws.Range("T5").Value = 0.0006058749567890857
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Il34"
$ws.Range("C6").Value = "Ptprz1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 2.952820666666666
$ws.Range("H6").Value = 8.858461999999999
$ws.Range("I6").Value = 0.09055216274676732
$ws.Range("J6").Value = 0.09055216274676733
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.1030763333333333
$ws.Range("N6").Value = 0.309229
$ws.Range("O6").Value = 0.01118975126488057
$ws.Range("P6").Value = 0.01118975126488057
$ws.Range("Q6").Value = 0.3043659273108889
$ws.Range("R6").Value = 2.739293345798
$ws.Range("S6").Value = 0.001013256177633311
$ws.Range("T6").Value = 0.001013256177633311
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Il34"
$ws.Range("C7").Value = "Ptprz1"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 2.952820666666666
$ws.Range("H7").Value = 8.858461999999999
$ws.Range("I7").Value = 0.09055216274676732
$ws.Range("J7").Value = 0.09055216274676733
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 9.046962666666667
$ws.Range("N7").Value = 27.140888
$ws.Range("O7").Value = 0.9821193543554519
$ws.Range("P7").Value = 0.9821193543554518
$ws.Range("Q7").Value = 26.71405833269511
$ws.Range("R7").Value = 240.426524994256
$ws.Range("S7").Value = 0.08893303161234493
$ws.Range("T7").Value = 0.08893303161234493
$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Il34"
$ws.Range("C8").Value = "Ptprz1"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 27.63137633333333
$ws.Range("H8").Value = 82.89412899999999
$ws.Range("I8").Value = 0.8473528090947983
$ws.Range("J8").Value = 0.8473528090947984
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.06163433333333333
$ws.Range("N8").Value = 0.184903
$ws.Range("O8").Value = 0.006690894379667537
$ws.Range("P8").Value = 0.006690894379667537
$ws.Range("Q8").Value = 1.703041459387444
$ws.Range("R8").Value = 15.327373134487
$ws.Range("S8").Value = 0.005669548147967885
$ws.Range("T8").Value = 0.005669548147967886
$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Il34"
$ws.Range("C9").Value = "Ptprz1"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 27.63137633333333
$ws.Range("H9").Value = 82.89412899999999
$ws.Range("I9").Value = 0.8473528090947983
$ws.Range("J9").Value = 0.8473528090947984
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.1030763333333333
$ws.Range("N9").Value = 0.309229
$ws.Range("O9").Value = 0.01118975126488057
$ws.Range("P9").Value = 0.01118975126488057
$ws.Range("Q9").Value = 2.848140957393444
$ws.Range("R9").Value = 25.633268616541
$ws.Range("S9").Value = 0.009481667167368627
$ws.Range("T9").Value = 0.009481667167368627
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Il34"
$ws.Range("C10").Value = "Ptprz1"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 27.63137633333333
$ws.Range("H10").Value = 82.89412899999999
$ws.Range("I10").Value = 0.8473528090947983
$ws.Range("J10").Value = 0.8473528090947984
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 9.046962666666667
$ws.Range("N10").Value = 27.140888
$ws.Range("O10").Value = 0.9821193543554519
$ws.Range("P10").Value = 0.9821193543554518
$ws.Range("Q10").Value = 249.9800301162836
$ws.Range("R10").Value = 2249.820271046552
$ws.Range("S10").Value = 0.8322015937794618
$ws.Range("T10").Value = 0.8322015937794618
Write-Output "Edit applied successfully"
